$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 05:53"

# Peru overtook Mexico: row 9 now holds Peru with fresh numbers, row 10 holds
# Mexico with the numbers that used to be in row 9 (unchanged).
$ws.Range("A9").Value = "Peru"
$ws.Range("B9").Value = 525803
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 359781
$ws.Range("E9").Value = 139947
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 26075

$ws.Range("A10").Value = "Mexico"
$ws.Range("B10").Value = 517714
$ws.Range("C10").Value = 6345
$ws.Range("D10").Value = 351372
$ws.Range("E10").Value = 109799
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 635
$ws.Range("H10").Value = 56543

# Kazajistan (row 29) data refresh
$ws.Range("B29").Value = 102696
$ws.Range("C29").Value = 409
$ws.Range("D29").Value = 81558
$ws.Range("E29").Value = 19869

# Venezuela (row 63) data refresh
$ws.Range("B63").Value = 32607
$ws.Range("D63").Value = 21747
$ws.Range("E63").Value = 10584
$ws.Range("H63").Value = 276

# Vietnam (row 159) data refresh
$ws.Range("E159").Value = 480
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 24

# Mongolia (row 180) data refresh
$ws.Range("D180").Value = 276
$ws.Range("E180").Value = 22

# Islas Malvinas overtook Montserrat: row 213 now holds Islas Malvinas with
# the numbers that used to be in row 214, row 214 holds Montserrat with the
# numbers that used to be in row 213.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
